$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 661
$ws.Range("F3").Value = 6490
$ws.Range("F4").Value = 1040
$ws.Range("F6").Value = 1436
$ws.Range("F7").Value = 3187
$ws.Range("F9").Value = 558
$ws.Range("F10").Value = 2104
$ws.Range("F11").Value = 450
$ws.Range("F12").Value = 379
$ws.Range("F13").Value = 221
$ws.Range("F14").Value = 109
$ws.Range("F15").Value = 241
$ws.Range("F16").Value = 1037
$ws.Range("F18").Value = 65
$ws.Range("F19").Value = 158
$ws.Range("F20").Value = 4070
$ws.Range("F21").Value = 1229
$ws.Range("F22").Value = 3190
$ws.Range("F24").Value = 105
$ws.Range("F25").Value = 2979
$ws.Range("F26").Value = 4629
$ws.Range("F27").Value = 117
$ws.Range("F29").Value = 510
$ws.Range("F30").Value = 3036
$ws.Range("F31").Value = 296
$ws.Range("F34").Value = 69
$ws.Range("F35").Value = 549
$ws.Range("F36").Value = 1096
$ws.Range("F37").Value = 1348
$ws.Range("F38").Value = 100
$ws.Range("F39").Value = 1223
$ws.Range("F40").Value = 792
$ws.Range("F42").Value = 729
$ws.Range("F43").Value = 476
$ws.Range("F44").Value = 41
$ws.Range("F45").Value = 205
$ws.Range("F46").Value = 38
$ws.Range("F47").Value = 86
$ws.Range("F48").Value = 347
$ws.Range("F49").Value = 3668

$ws = $wb.Worksheets.Item(2)
$ws.Range("F11").Value = 5

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 1692

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 6490
$ws.Range("F4").Value = 1693
$ws.Range("F6").Value = 1436
$ws.Range("F7").Value = 3187
$ws.Range("F9").Value = 2104
$ws.Range("F10").Value = 450
$ws.Range("F11").Value = 379
$ws.Range("F13").Value = 221
$ws.Range("F16").Value = 109
$ws.Range("F17").Value = 241
$ws.Range("F18").Value = 1037
$ws.Range("F19").Value = 5
$ws.Range("F21").Value = 158
$ws.Range("F22").Value = 4070
$ws.Range("F24").Value = 1229
$ws.Range("F26").Value = 3190
$ws.Range("F27").Value = 2979
$ws.Range("F28").Value = 4629
$ws.Range("F30").Value = 3036
$ws.Range("F31").Value = 296
$ws.Range("F34").Value = 549
$ws.Range("F35").Value = 1096
$ws.Range("F36").Value = 1348
$ws.Range("F37").Value = 100
$ws.Range("F38").Value = 1223
$ws.Range("F39").Value = 792
$ws.Range("F41").Value = 476
$ws.Range("F43").Value = 41
$ws.Range("F45").Value = 205
$ws.Range("F46").Value = 38
$ws.Range("F47").Value = 86
$ws.Range("F48").Value = 347
$ws.Range("F49").Value = 3668
